$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "SamplesTab" query (row 3 / B3): dropped the Tumor and Analyte Type
# columns from the SELECT list (this is the "All studies" sample query).
$newSampleQuery = @'
SELECT
    DISTINCT (smp.sample_id) AS "Sample ID",
    sp.participant_id AS "Participant ID", 
    s.study_name AS "Study Name",
    s.phs_accession AS Accession
FROM 
    df_participant sp
JOIN 
    df_study s ON sp."study.phs_accession" = s.phs_accession
JOIN 
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
JOIN
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
WHERE 
   s.phs_accession = 'phs001819' AND gi.reference_genome_assembly = 'GRCh37'
ORDER BY 
    smp.sample_id ASC
LIMIT 100;
'@

$ws.Range("B3").Value = $newSampleQuery

# Move the selection/view: top-left visible row now row 3, active cell B3.
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 3
    $win.ScrollColumn = 1
} catch {
    # View-scroll state isn't critical; ignore if unsupported by the host.
}
$ws.Range("B3").Select()
